$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '88.591.14'
$ws.Range("E2").Value = '  +3.60%  '

Set-TextValue $ws.Range("D3") '3.191.40'
$ws.Range("E3").Value = '  -1.22%  '

Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.14%  '

Set-TextValue $ws.Range("D5") '210.38'
$ws.Range("E5").Value = '  +1.34%  '

Set-TextValue $ws.Range("D6") '615.89'
$ws.Range("E6").Value = '  -0.69%  '

Set-TextValue $ws.Range("D7") '0.391'
$ws.Range("E7").Value = '  +9.32%  '

Set-TextValue $ws.Range("D8") '0.692'
$ws.Range("E8").Value = '  +7.33%  '

Set-TextValue $ws.Range("D9") '0.999'
$ws.Range("E9").Value = '  -0.18%  '

Set-TextValue $ws.Range("D10") '3.184.69'
$ws.Range("E10").Value = '  -1.53%  '

Set-TextValue $ws.Range("D11") '0.551'
$ws.Range("E11").Value = '  -3.38%  '

$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("E13").Value = '  +0.08%  '

Set-TextValue $ws.Range("D14") '5.34'
$ws.Range("E14").Value = '  +2.24%  '

Set-TextValue $ws.Range("D15") '3.775.21'
$ws.Range("E15").Value = '  -1.76%  '

Set-TextValue $ws.Range("D16") '32.84'
$ws.Range("E16").Value = '  -1.38%  '

Set-TextValue $ws.Range("D17") '88.269.79'
$ws.Range("E17").Value = '  +3.00%  '

Set-TextValue $ws.Range("D18") '3.188.73'
$ws.Range("E18").Value = '  -1.69%  '

Set-TextValue $ws.Range("D19") '3.22'
$ws.Range("E19").Value = '  +9.26%  '

Set-TextValue $ws.Range("D20") '13.64'
$ws.Range("E20").Value = '  -1.64%  '

Set-TextValue $ws.Range("D21") '418.25'
$ws.Range("E21").Value = '  -1.42%  '

Set-TextValue $ws.Range("D22") '8.57'
$ws.Range("E22").Value = '  -3.30%  '

Set-TextValue $ws.Range("D23") '5.17'
$ws.Range("E23").Value = '  -1.48%  '

Set-TextValue $ws.Range("D24") '0.0000158'
$ws.Range("E24").Value = '  +24.57%  '

Set-TextValue $ws.Range("D25") '5.37'
$ws.Range("E25").Value = '  +5.20%  '

Set-TextValue $ws.Range("D26") '12.48'
$ws.Range("E26").Value = '  +2.02%  '

$ws.Range("E27").Value = '  -2.87%  '

Set-TextValue $ws.Range("D28") '74.07'
$ws.Range("E28").Value = '  -1.90%  '

Set-TextValue $ws.Range("D29") '0.998'
$ws.Range("E29").Value = '  +0.21%  '

Set-TextValue $ws.Range("D30") '0.169'
$ws.Range("E30").Value = '  -1.30%  '

Set-TextValue $ws.Range("D31") '1.00'
$ws.Range("E31").Value = '  +0.15%  '

Set-TextValue $ws.Range("D32") '556.09'
$ws.Range("E32").Value = '  +2.73%  '

Set-TextValue $ws.Range("D33") '8.37'
$ws.Range("E33").Value = '  -4.11%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D34") '7.01'
$ws.Range("E34").Value = '  +4.27%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D35") '1.33'
$ws.Range("E35").Value = '  -5.83%  '

Set-TextValue $ws.Range("D36") '1.88'
$ws.Range("E36").Value = '  -3.07%  '

$ws.Range("E37").Value = '  -2.25%  '

Set-TextValue $ws.Range("D38") '22.14'
$ws.Range("E38").Value = '  -0.20%  '

$ws.Range("E39").Value = '  +0.88%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D40") '0.997'
$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D41") '3.18'
$ws.Range("E41").Value = '  +10.52%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D42") '1.95'
$ws.Range("E42").Value = '  -1.00%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D43") '1.00'
$ws.Range("E43").Value = '  +0.07%  '

Set-TextValue $ws.Range("D44") '0.380'
$ws.Range("E44").Value = '  -2.24%  '

Set-TextValue $ws.Range("D45") '150.42'
$ws.Range("E45").Value = '  -4.65%  '

Set-TextValue $ws.Range("D46") '176.12'
$ws.Range("E46").Value = '  +0.55%  '

Set-TextValue $ws.Range("D47") '43.30'
$ws.Range("E47").Value = '  -1.18%  '

Set-TextValue $ws.Range("D48") '0.127'
$ws.Range("E48").Value = '  +8.10%  '

Set-TextValue $ws.Range("D49") '1.26'
$ws.Range("E49").Value = '  -3.77%  '

Set-TextValue $ws.Range("D50") '24.62'
$ws.Range("E50").Value = '  +3.33%  '

$ws.Range("E51").Value = '  -4.87%  '
